# Swap the contents of columns D ("codeforiati:group-name") and E
# ("codeforiati:group-code"), including the header row, for the full
# used range of the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $temp = $dCell.Value()
    $dCell.Value = $eCell.Value()
    $eCell.Value = $temp
}
